$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F
$ws.Range("F1").Value = "sd_temp"

# New sd_temp values for rows 2-12
$values = @(
    3.2247637885475,
    2.77024475400633,
    1.66322373635078,
    3.09716239692286,
    2.43240550981441,
    1.83864760081969,
    2.81217129758559,
    2.43730365892953,
    1.12006493318265,
    3.39912052233226,
    2.90311680708871
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
